$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: iaest-measure: -> iaest-dimension: for horas-trabajadas / situacion-profesional
$ws.Range("A3").Value = "iaest-dimension:horas-trabajadas"
$ws.Range("E3").Value = "iaest-dimension:situacion-profesional"

# Row 4: these columns become dimensions ("dim") instead of measures ("medida")
$ws.Range("A4").Value = "dim"
$ws.Range("E4").Value = "dim"

# Row 5: data type changes from xsd:string to skos:Concept
$ws.Range("A5").Value = "skos:Concept"
$ws.Range("E5").Value = "skos:Concept"

# Row 6 (new): mapping file references
$ws.Range("A6").Value = "mapping-horas-trabajadas.xlsx"
$ws.Range("E6").Value = "mapping-situacion-profesional.xlsx"

# Match the styling of the rest of the sheet (same font/style as row 5)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("E5").Copy()
$ws.Range("E6").PasteSpecial(-4122)
